$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "Flip Zombie Image"
$ws.Range("H7").Interior.Color = $ws.Range("G4").Interior.Color
$ws.Range("H7").HorizontalAlignment = $ws.Range("G4").HorizontalAlignment
$ws.Range("H7").WrapText = $ws.Range("G4").WrapText
$ws.Range("H7").Select()
